$d = $word.ActiveDocument

# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# MatchCase=$true (avoid cross-matching differently-cased headings) and
# MatchWholeWord=$false (avoid disturbing/merging adjacent empty runs).

# Title text updated throughout (heading, and bold run near the end)
$d.Content.Find.Execute("Play Blazing Goddess for Free | Review of Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Blazing Goddess for Free", 2) | Out-Null

# "What we like" bullet list updates
$d.Content.Find.Execute("Beautiful graphics and design", $true, $false, $false, $false, $false, $true, 1, $false, "Beautiful graphics and theme", 2) | Out-Null
$d.Content.Find.Execute("Extra wilds increase chances of winning", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting free spins and bonus rounds", 2) | Out-Null
$d.Content.Find.Execute("Free spins and bonus rounds", $true, $false, $false, $false, $false, $true, 1, $false, "Playable on all devices", 2) | Out-Null

# "What we don't like" bullet list update
$d.Content.Find.Execute("Not suitable for those looking for huge payouts", $true, $false, $false, $false, $false, $true, 1, $false, "Limited range of features", 2) | Out-Null

# Meta description text near the end
$d.Content.Find.Execute("Find out all about the gameplay mechanics, graphics, symbols, and ways to win in Blazing Goddess. Play for free and read the game review now.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Blazing Goddess slot game and play it for free on all devices.", 2) | Out-Null
